# Helper to build a VBA-style packed RGB integer (0x00BBGGRR) from R,G,B bytes,
# matching the PowerPoint COM ForeColor.RGB convention used by this host.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$C00000 = RGB 0xC0 0x00 0x00   # srgbClr C00000 (red) - "significant" highlight color
$AFABAB = RGB 0xAF 0xAB 0xAB   # flattened bg2 / lumMod 75% (theme "gray") - "non-significant" color
$ACCENT1 = RGB 0x44 0x72 0xC4  # flattened accent1 theme color (unused target, kept for reference)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Connectors / boxes whose line weight only changes (no color/dash/text change) ---
$s.Shapes.Item("Straight Arrow Connector 72").Line.Weight = 9.51    # 118491 -> 120777
$s.Shapes.Item("Straight Arrow Connector 42").Line.Weight = 2.87    # 35306  -> 36449
$s.Shapes.Item("Straight Arrow Connector 43").Line.Weight = 3.79    # 38862  -> 48133
$s.Shapes.Item("Straight Arrow Connector 49").Line.Weight = 3.31    # 34036  -> 42037
$s.Shapes.Item("Straight Arrow Connector 55").Line.Weight = 2.79    # 33147  -> 35433
$s.Shapes.Item("Straight Arrow Connector 57").Line.Weight = 2.9     # 33528  -> 36830
$s.Shapes.Item("Straight Arrow Connector 59").Line.Weight = 3.21    # 38608  -> 40767
$s.Shapes.Item("TextBox 61").Line.Weight = 2.9                      # 33528  -> 36830
$s.Shapes.Item("Straight Arrow Connector 63").Line.Weight = 2.81    # 33909  -> 35687
$s.Shapes.Item("Straight Arrow Connector 76").Line.Weight = 7.27    # 76454  -> 92329
$s.Shapes.Item("Straight Arrow Connector 1").Line.Weight = 5.84     # 65659  -> 74168
$s.Shapes.Item("Straight Arrow Connector 29").Line.Weight = 5.46    # 65659  -> 69342

# --- Connectors that flip from dashed/gray to solid/red (weight + color + dash) ---
$cxn53 = $s.Shapes.Item("Straight Arrow Connector 53")
$cxn53.Line.Weight = 3.19            # 31623 -> 40513
$cxn53.Line.ForeColor.RGB = $C00000  # bg2/lumMod75% -> C00000
$cxn53.Line.DashStyle = 1            # dash -> solid

$cxn3 = $s.Shapes.Item("Straight Arrow Connector 3")
$cxn3.Line.Weight = 3.04             # 62102 -> 38608
$cxn3.Line.ForeColor.RGB = $AFABAB   # accent1 -> bg2/lumMod75%
$cxn3.Line.DashStyle = 4             # solid -> dash

# --- Value-label textboxes: weight (+ color/dash where relevant) + displayed number ---

# -0.079 -> -0.109, plus the label nudges down slightly (Top changes)
$tb60 = $s.Shapes.Item("TextBox 60")
$tb60.Line.Weight = 3.21             # 38608 -> 40767
$tb60.Top = 375.7185039370079        # off y 4760608 -> 4771625
$tb60.TextFrame.TextRange.Text = "-0.109"

# 0.026 -> 0.036
$tb64 = $s.Shapes.Item("TextBox 64")
$tb64.Line.Weight = 2.81             # 33909 -> 35687
$tb64.TextFrame.TextRange.Text = "0.036"

# -0.059 -> -0.127
$tb65 = $s.Shapes.Item("TextBox 65")
$tb65.Line.Weight = 3.31             # 33655 -> 42037
$tb65.TextFrame.TextRange.Text = "-0.127"

# -0.039 -> -0.104, line flips dashed/gray -> solid/red, and the run becomes bold
$tb66 = $s.Shapes.Item("TextBox 66")
$tb66.Line.Weight = 3.19             # 31623 -> 40513
$tb66.Line.ForeColor.RGB = $C00000   # bg2/lumMod75% -> C00000
$tb66.Line.DashStyle = 1             # dash -> solid
$tb66.TextFrame.TextRange.Font.Bold = $true
$tb66.TextFrame.TextRange.Text = "-0.104"

# 0.020 -> 0.032
$tb67 = $s.Shapes.Item("TextBox 67")
$tb67.Line.Weight = 2.79             # 33147 -> 35433
$tb67.TextFrame.TextRange.Text = "0.032"

# -0.173 -> -0.214
$tb68 = $s.Shapes.Item("TextBox 68")
$tb68.Line.Weight = 3.79             # 38862 -> 48133
$tb68.TextFrame.TextRange.Text = "-0.214"

# -0.081 -> -0.046
$tb70 = $s.Shapes.Item("TextBox 70")
$tb70.Line.Weight = 2.87             # 35306 -> 36449
$tb70.TextFrame.TextRange.Text = "-0.046"

# -0.902 -> -1.254
$tb73 = $s.Shapes.Item("TextBox 73")
$tb73.Line.Weight = 9.51             # 118491 -> 120777
$tb73.TextFrame.TextRange.Text = "-1.254"

# -0.774 -> -0.848
$tb79 = $s.Shapes.Item("TextBox 79")
$tb79.Line.Weight = 7.27             # 76454 -> 92329
$tb79.TextFrame.TextRange.Text = "-0.848"

# -0.602 -> -0.587
$tb2 = $s.Shapes.Item("TextBox 2")
$tb2.Line.Weight = 5.84              # 65659 -> 74168
$tb2.TextFrame.TextRange.Text = "-0.587"

# -0.540 -> 0.078, line flips solid/accent1 -> dashed/gray, box narrows, run unbolds
$tb23 = $s.Shapes.Item("TextBox 23")
$tb23.Line.Weight = 3.04             # 62102 -> 38608
$tb23.Line.ForeColor.RGB = $AFABAB   # accent1 -> bg2/lumMod75%
$tb23.Line.DashStyle = 4             # solid -> dash
$tb23.Width = 69.06795275590551      # ext cx 979755 -> 877163
$tb23.TextFrame.TextRange.Font.Bold = $false
$tb23.TextFrame.TextRange.Text = "0.078"

# 0.527 -> 0.517
$tb74 = $s.Shapes.Item("TextBox 74")
$tb74.Line.Weight = 5.46             # 65659 -> 69342
$tb74.TextFrame.TextRange.Text = "0.517"
